# fix createdAt mockup data
# Adds "createdAt" (H) and "updatedAt" (I) columns to the patient mockup
# sheet. Both columns carry the same timestamp per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("H1").Value = "createdAt"
$ws.Range("I1").Value = "updatedAt"

$timestamps = @(
    "2023-04-16T03:42:46.836Z",
    "2023-04-09T03:20:31.077Z",
    "2023-04-19T12:32:13.787Z",
    "2023-04-20T18:59:51.241Z",
    "2023-04-12T20:57:30.058Z",
    "2023-04-18T13:51:56.194Z",
    "2023-04-20T02:09:40.846Z",
    "2023-04-21T01:45:38.618Z",
    "2023-04-23T04:44:42.471Z",
    "2023-05-02T13:49:47.409Z",
    "2023-05-02T10:13:51.564Z",
    "2023-05-03T12:04:43.165Z",
    "2023-05-02T12:12:13.632Z",
    "2023-05-16T18:31:08.185Z",
    "2023-05-08T08:54:59.583Z",
    "2023-05-07T11:40:13.109Z",
    "2023-05-08T12:27:27.314Z",
    "2023-05-10T21:47:27.523Z",
    "2023-05-09T04:25:04.137Z"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ts = $timestamps[$i]
    $ws.Cells.Item($row, 8).Value = $ts
    $ws.Cells.Item($row, 9).Value = $ts
}

# Match the authored column widths for the two new columns.
$ws.Range("H1:I1").ColumnWidth = 23.21875

# Mirror the selection move captured in the diff.
$ws.Range("E6").Select()
